$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge-weight metrics (Tgfb3-Tgfbr2) per Dr Hou advice.
# Columns updated for every data row (2-17): E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$cols = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)

$data = @{
    2 = @(3,1.277420333333333,3.832261,0.01913942624337554,0.01913942624337554,3,28.21795633333333,84.653869,0.2124568395711989,0.212456839571199,36.04619118531211,324.415720667809,0.004066302010873632,0.004066302010873633)
    3 = @(3,1.277420333333333,3.832261,0.01913942624337554,0.01913942624337554,3,63.95730733333333,191.871922,0.4815432848151522,0.4815432848151524,81.70036485284911,735.303283675642,0.009216462182712389,0.00921646218271239)
    4 = @(3,1.277420333333333,3.832261,0.01913942624337554,0.01913942624337554,3,27.47719633333334,82.431589,0.2068795565595709,0.2068795565595709,35.09992929919211,315.899363692729,0.003959556014034147,0.003959556014034147)
    5 = @(3,1.277420333333333,3.832261,0.01913942624337554,0.01913942624337554,3,13.164899,39.494697,0.09912031905407785,0.09912031905407789,16.81710966887967,151.353987019917,0.001897106035755375,0.001897106035755375)
    6 = @(3,45.44725166666667,136.341755,0.6809303864519871,0.6809303864519872,3,28.21795633333333,84.653869,0.2124568395711989,0.212456839571199,1282.428563000011,11541.8570670001,0.1446683178735843,0.1446683178735844)
    7 = @(3,45.44725166666667,136.341755,0.6809303864519871,0.6809303864519872,3,63.95730733333333,191.871922,0.4815432848151522,0.4815432848151524,2906.683842300345,26160.15458070311,0.3278974550225409,0.327897455022541)
    8 = @(3,45.44725166666667,136.341755,0.6809303864519871,0.6809303864519872,3,27.47719633333334,82.431589,0.2068795565595709,0.2068795565595709,1248.763056855411,11238.8675116987,0.1408705763971244,0.1408705763971244)
    9 = @(3,45.44725166666667,136.341755,0.6809303864519871,0.6809303864519872,3,13.164899,39.494697,0.09912031905407785,0.09912031905407789,598.3084780192484,5384.776302173235,0.06749403715873749,0.06749403715873753)
    10 = @(3,0.4966396666666666,1.489919,0.007441088905245192,0.007441088905245193,3,28.21795633333333,84.653869,0.2124568395711989,0.212456839571199,14.01415642740122,126.127407846611,0.001580910231776706,0.001580910231776707)
    11 = @(3,0.4966396666666666,1.489919,0.007441088905245192,0.007441088905245193,3,63.95730733333333,191.871922,0.4815432848151522,0.4815432848151524,31.76373579492422,285.873622154318,0.003583206394033355,0.003583206394033356)
    12 = @(3,0.4966396666666666,1.489919,0.007441088905245192,0.007441088905245193,3,27.47719633333334,82.431589,0.2068795565595709,0.2068795565595709,13.64626562792122,122.816390651291,0.001539409173037468,0.001539409173037469)
    13 = @(3,0.4966396666666666,1.489919,0.007441088905245192,0.007441088905245193,3,13.164899,39.494697,0.09912031905407785,0.09912031905407789,6.538211051060333,58.843899459543,0.0007375631063976623,0.0007375631063976627)
    14 = @(3,19.52156333333333,58.56469,0.2924890983993922,0.2924890983993922,3,28.21795633333333,84.653869,0.2124568395711989,0.212456839571199,550.8586216984011,4957.72759528561,0.06214130945496427,0.0621413094549643)
    15 = @(3,19.52156333333333,58.56469,0.2924890983993922,0.2924890983993922,3,63.95730733333333,191.871922,0.4815432848151522,0.4815432848151524,1248.546625737131,11236.91963163418,0.1408461612158656,0.1408461612158657)
    16 = @(3,19.52156333333333,58.56469,0.2924890983993922,0.2924890983993922,3,27.47719633333334,82.431589,0.2068795565595709,0.2068795565595709,536.3978284436012,4827.58045599241,0.06051001497537496,0.06051001497537498)
    17 = @(3,19.52156333333333,58.56469,0.2924890983993922,0.2924890983993922,3,13.164899,39.494697,0.09912031905407785,0.09912031905407789,256.9994096054367,2312.99468644893,0.02899161275318732,0.02899161275318734)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}